# Insert a new data row for "Rabanito" (Mercado Mayorista Lo Valledor de
# Santiago) above the existing row 25, shifting the subsequent rows (old
# 25-30) down to (26-31) and growing the used range to A1:R31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before row 25; everything below (including formatting)
# shifts down by one.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new record.
$ws.Cells.Item(25, 1).Value = 6
$ws.Cells.Item(25, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(25, 3).Value = "Metropolitana"
$ws.Cells.Item(25, 4).Value = 44855
$ws.Cells.Item(25, 5).Value = 13
$ws.Cells.Item(25, 6).Value = 300000001
$ws.Cells.Item(25, 7).Value = "Rabanito"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 7900
$ws.Cells.Item(25, 11).Value = 3000
$ws.Cells.Item(25, 12).Value = 3000
$ws.Cells.Item(25, 13).Value = 3000
$ws.Cells.Item(25, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(25, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(25, 16).Value = 30
$ws.Cells.Item(25, 17).Value = 100
$ws.Cells.Item(25, 18).Value = "Hortaliza"
